# Czech Republic CFL Group B 2023-2024 — apply upstream data refresh.
# Net effect per the upload diff:
#   * rows 5 & 6 swap their match data (cols F:V)
#   * rows 35-38 rotate their match data up by one (35<-36<-37<-38<-35)
#   * rows 44 & 45 swap their match data (cols F:V)
#   * three brand-new match rows are appended at the end (58, 59, 60)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F..V (6..22) hold the per-match payload; A..E (index/pais/torneio/
# temporada/data_partida) are untouched by every one of these edits.
$cols = @(6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)

function Swap-RowPayload($rowA, $rowB, $cols) {
    foreach ($c in $cols) {
        $va = $ws.Cells.Item($rowA, $c).Value()
        $vb = $ws.Cells.Item($rowB, $c).Value()
        $ws.Cells.Item($rowA, $c).Value = $vb
        $ws.Cells.Item($rowB, $c).Value = $va
    }
}

function Rotate-RowPayloadUp($rowList, $cols) {
    # destination row i receives the payload that used to live in rowList[i+1]
    # (wrapping around), i.e. everything shifts "up" by one row.
    $old = @{}
    foreach ($r in $rowList) {
        foreach ($c in $cols) {
            $old["$r,$c"] = $ws.Cells.Item($r, $c).Value()
        }
    }
    for ($i = 0; $i -lt $rowList.Length; $i++) {
        $destRow = $rowList[$i]
        $srcRow = $rowList[($i + 1) % $rowList.Length]
        foreach ($c in $cols) {
            $ws.Cells.Item($destRow, $c).Value = $old["$srcRow,$c"]
        }
    }
}

# --- rows 5 & 6 swap ---
Swap-RowPayload 5 6 $cols

# --- rows 35-38 rotate ---
Rotate-RowPayloadUp @(35,36,37,38) $cols

# --- rows 44 & 45 swap ---
Swap-RowPayload 44 45 $cols

# --- append three new rows (58, 59, 60) ---
# Copy number/text formatting from the last existing data row (57) so the
# new "Indice" (A) and "data_partida" (E) cells keep the same styles as
# every other row, then fill in the values.
$ws.Range("A57").Copy()
$ws.Range("A58:A60").PasteSpecial(-4122)
$ws.Range("E57").Copy()
$ws.Range("E58:E60").PasteSpecial(-4122)

$newRows = @(
    @{ Row=58; A=57; E=45191.70833333334; F="Kolin"; G=6; H="Prepere"; I=1;
       J=1.94; K="21/09/2023 04:42"; L=2.05; M="22/09/2023 15:16";
       N=3.56; O="21/09/2023 04:42"; P=3.71; Q="22/09/2023 16:02";
       R=3;    S="21/09/2023 04:42"; T=3.05; U="22/09/2023 15:15";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/kolin-prepere/UgFWi9bU/" },
    @{ Row=59; A=58; E=45191.75; F="Usti nad Labem"; G=3; H="Teplice B"; I=0;
       J=1.34; K="21/09/2023 05:12"; L=1.3;  M="22/09/2023 17:50";
       N=4.85; O="21/09/2023 05:12"; P=5.96; Q="22/09/2023 17:51";
       R=5.96; S="21/09/2023 05:12"; T=6.79; U="22/09/2023 17:51";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/usti-nad-labem-teplice/fHtwDhUu/" },
    @{ Row=60; A=59; E=45192.42708333334; F="Jablonec B"; G=2; H="Liberec B"; I=2;
       J=1.65; K="21/09/2023 21:42"; L=1.76; M="23/09/2023 10:00";
       N=3.86; O="21/09/2023 21:42"; P=4.35; Q="23/09/2023 10:08";
       R=3.8;  S="21/09/2023 21:42"; T=3.47; U="23/09/2023 10:08";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/jablonec-liberec/bebLfBEB/" }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value  = $rowData.A
    $ws.Cells.Item($r, 2).Value  = "czech-republic"
    $ws.Cells.Item($r, 3).Value  = "cfl-group-b"
    $ws.Cells.Item($r, 4).Value  = "2023-2024"
    $ws.Cells.Item($r, 5).Value  = $rowData.E
    $ws.Cells.Item($r, 6).Value  = $rowData.F
    $ws.Cells.Item($r, 7).Value  = $rowData.G
    $ws.Cells.Item($r, 8).Value  = $rowData.H
    $ws.Cells.Item($r, 9).Value  = $rowData.I
    $ws.Cells.Item($r, 10).Value = $rowData.J
    $ws.Cells.Item($r, 11).Value = $rowData.K
    $ws.Cells.Item($r, 12).Value = $rowData.L
    $ws.Cells.Item($r, 13).Value = $rowData.M
    $ws.Cells.Item($r, 14).Value = $rowData.N
    $ws.Cells.Item($r, 15).Value = $rowData.O
    $ws.Cells.Item($r, 16).Value = $rowData.P
    $ws.Cells.Item($r, 17).Value = $rowData.Q
    $ws.Cells.Item($r, 18).Value = $rowData.R
    $ws.Cells.Item($r, 19).Value = $rowData.S
    $ws.Cells.Item($r, 20).Value = $rowData.T
    $ws.Cells.Item($r, 21).Value = $rowData.U
    $ws.Cells.Item($r, 22).Value = $rowData.V
}
